$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.396.89'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '''1.893.68'
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("D4").Value = '''1.014'
$ws.Range("E4").Value = '  +0.99%  '
$ws.Range("D5").Value = '''316.83'
$ws.Range("E5").Value = '  +1.48%  '
$ws.Range("D6").Value = '''1.013'
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("D7").Value = '''0.5172'
$ws.Range("E7").Value = '  +1.62%  '
$ws.Range("D8").Value = '''0.3927'
$ws.Range("E8").Value = '  +2.10%  '
$ws.Range("D9").Value = '''0.08429'
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("D10").Value = '''1.129'
$ws.Range("E10").Value = '  +1.44%  '
$ws.Range("D11").Value = '''42.02'
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '''1.916.30'
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''6.290'
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").Value = '''20.73'
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").Value = '''7.322'
$ws.Range("E15").Value = '  +1.43%  '
$ws.Range("D16").Value = '''1.015'
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = '''91.60'
$ws.Range("E18").Value = '  +0.81%  '
$ws.Range("D19").Value = '''0.06745'
$ws.Range("E19").Value = '  +1.74%  '
$ws.Range("D20").Value = '''17.90'
$ws.Range("E20").Value = '  +1.26%  '
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("D22").Value = '''6.078'
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("D23").Value = '''28.470.44'
$ws.Range("E23").Value = '  +1.50%  '
$ws.Range("D24").Value = '''11.21'
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("D25").Value = '''2.268'
$ws.Range("E25").Value = '  +1.86%  '
$ws.Range("D26").Value = '''160.73'
$ws.Range("E26").Value = '  +1.70%  '
$ws.Range("D27").Value = '''2.486'
$ws.Range("D28").Value = '''20.79'
$ws.Range("E28").Value = '  +1.26%  '
$ws.Range("D29").Value = '''127.09'
$ws.Range("E29").Value = '  +1.85%  '
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("D31").Value = '''1.043'
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("D32").Value = '''5.842'
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("D33").Value = '''3.636'
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D34").Value = '''9.686'
$ws.Range("E34").Value = '  +2.66%  '
$ws.Range("D35").Value = '''0.02473'
$ws.Range("E35").Value = '  +2.27%  '
$ws.Range("D36").Value = '''0.06617'
$ws.Range("E36").Value = '  +1.43%  '
$ws.Range("D37").Value = '''0.2219'
$ws.Range("E37").Value = '  +2.23%  '
$ws.Range("D38").Value = '''1.207'
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("D39").Value = '''0.6532'
$ws.Range("E39").Value = '  +1.03%  '
$ws.Range("D40").Value = '''1.250'
$ws.Range("E40").Value = '  +2.07%  '
$ws.Range("D41").Value = '''5.014'
$ws.Range("E41").Value = '  +1.02%  '
$ws.Range("D42").Value = '''11.36'
$ws.Range("E42").Value = '  +1.45%  '
$ws.Range("D43").Value = '''0.6146'
$ws.Range("E43").Value = '  +0.86%  '
$ws.Range("D44").Value = '''13.24'
$ws.Range("E44").Value = '  +1.00%  '
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("D46").Value = '''1.287'
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").Value = '''2.029'
$ws.Range("E47").Value = '  +0.57%  '
$ws.Range("D48").Value = '''1.243'
$ws.Range("E48").Value = '  +2.91%  '
$ws.Range("D49").Value = '''121.46'
$ws.Range("E49").Value = '  +1.19%  '
$ws.Range("D50").Value = '''0.06947'
$ws.Range("E50").Value = '  +1.44%  '
$ws.Range("D51").Value = '''78.32'
$ws.Range("E51").Value = '  +0.05%  '
